$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("SIU")

# Row 4 - Heydi Quemé
$ws.Range("C4").Value = 8
$ws.Range("H4").Value = 10
$ws.Range("M4").Value = 10
$ws.Range("R4").Value = 9
$ws.Range("W4").Value = 9

# Row 5 - Brayan Cifuentes
$ws.Range("C5").Value = 10
$ws.Range("H5").Value = 10
$ws.Range("M5").Value = 10
$ws.Range("R5").Value = 10
$ws.Range("W5").Value = 10

# Row 6 - Sebastian Moreira
$ws.Range("C6").Value = 10
$ws.Range("H6").Value = 10
$ws.Range("M6").Value = 10
$ws.Range("R6").Value = 10
$ws.Range("W6").Value = 10

# Update selection to W5 as recorded in the saved workbook state
$ws.Range("W5").Select()

$wb.Save()
